$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 4
$ws.Range("D4").Value = 2
$ws.Range("F4").Value = -3
$ws.Range("H4").Value = 46

# Update the active selection to E4
$null = $ws.Range("E4").Select()
